# Fruta / hortaliza, semanal
# Insert a new weekly record at row 158 of the "Ají" sheet (Feria Lagunitas
# de Puerto Montt), pushing the existing rows 158-167 down to 159-168.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 158, shifting rows 158:167
# down to 159:168 (Excel copies the formatting - including the date
# number-format on column D - from the row above, same as interactive use).
$ws.Rows.Item(158).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the newly inserted row 158 with the new weekly data point.
$ws.Cells.Item(158, 1).Value  = 4
$ws.Cells.Item(158, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(158, 3).Value  = "Los Lagos"
$ws.Cells.Item(158, 4).Value  = 44516
$ws.Cells.Item(158, 5).Value  = 10
$ws.Cells.Item(158, 6).Value  = 100112021
$ws.Cells.Item(158, 7).Value  = "Ají"
$ws.Cells.Item(158, 8).Value  = "Inferno"
$ws.Cells.Item(158, 9).Value  = "Primera"
$ws.Cells.Item(158, 10).Value = 160
$ws.Cells.Item(158, 11).Value = 26000
$ws.Cells.Item(158, 12).Value = 27000
$ws.Cells.Item(158, 13).Value = 26500
$ws.Cells.Item(158, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(158, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(158, 16).Value = 2208
$ws.Cells.Item(158, 17).Value = 12
$ws.Cells.Item(158, 18).Value = "Hortaliza"
